$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4724.625
$ws.Range("I18").Value = 6099.5
$ws.Range("K18").Value = 6099.5
$ws.Range("M18").Value = -5815.5
$ws.Range("H43").Value = 362.0625
$ws.Range("I43").Value = 100.888885
$ws.Range("J43").Value = 697.8570999999999
$ws.Range("K43").Value = 100.888885
$ws.Range("L43").Value = 697.8570999999999
$ws.Range("M43").Value = -31.888885
$ws.Range("N43").Value = -835.8570999999999
$ws.Range("H92").Value = 480.875
$ws.Range("I92").Value = 472
$ws.Range("K92").Value = 472
$ws.Range("M92").Value = 776
$ws.Range("H129").Value = 776.5
$ws.Range("I129").Value = 532
$ws.Range("J129").Value = 1999
$ws.Range("K129").Value = 1596
$ws.Range("L129").Value = 5997
$ws.Range("M129").Value = 3404
$ws.Range("N129").Value = -15997
$ws.Range("H132").Value = 1204.0444
$ws.Range("I132").Value = 1128.1666
$ws.Range("K132").Value = 3384.4998
$ws.Range("M132").Value = -854.4998000000001
$ws.Range("H137").Value = 2697.8064
$ws.Range("I137").Value = 2552.8572
$ws.Range("J137").Value = 2817.1765
$ws.Range("K137").Value = 7658.571599999999
$ws.Range("L137").Value = 8451.529500000001
$ws.Range("M137").Value = -5108.571599999999
$ws.Range("N137").Value = -13551.5295
$ws.Range("H141").Value = 2569.6667
$ws.Range("I141").Value = 2451
$ws.Range("J141").Value = 3074
$ws.Range("K141").Value = 7353
$ws.Range("L141").Value = 9222
$ws.Range("M141").Value = -2173
$ws.Range("N141").Value = -19582

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2714.6667
$ws.Range("I2").Value = 1214.2667
$ws.Range("K2").Value = 1214.2667
$ws.Range("M2").Value = -1101.2667
$ws.Range("H32").Value = 1670926.5
$ws.Range("I32").Value = 1869349.2
$ws.Range("J32").Value = 9136.375
$ws.Range("K32").Value = 1869349.2
$ws.Range("L32").Value = 9136.375
$ws.Range("M32").Value = -1869062.2
$ws.Range("N32").Value = -9710.375
$ws.Range("H45").Value = 6734.263
$ws.Range("I45").Value = 1862.7778
$ws.Range("K45").Value = 1862.7778
$ws.Range("M45").Value = -1485.7778
$ws.Range("H110").Value = 13890163
$ws.Range("I110").Value = 1160.4706
$ws.Range("K110").Value = 1160.4706
$ws.Range("M110").Value = 884.5293999999999
$ws.Range("H116").Value = 2714.6667
$ws.Range("I116").Value = 1214.2667
$ws.Range("K116").Value = 1214.2667
$ws.Range("M116").Value = 1079.7333
$ws.Range("H140").Value = 45643
$ws.Range("J140").Value = 45643
$ws.Range("L140").Value = 45643
$ws.Range("N140").Value = -56003

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2714.6667
$ws.Range("I3").Value = 1214.2667
$ws.Range("K3").Value = 1214.2667
$ws.Range("M3").Value = -1100.2667
$ws.Range("H20").Value = 8335350
$ws.Range("I20").Value = 12822224
$ws.Range("J20").Value = 2584.8572
$ws.Range("K20").Value = 12822224
$ws.Range("L20").Value = 2584.8572
$ws.Range("M20").Value = -12821977
$ws.Range("N20").Value = -3078.8572
$ws.Range("H22").Value = 256.57144
$ws.Range("I22").Value = 256.57144
$ws.Range("K22").Value = 256.57144
$ws.Range("M22").Value = -83.57144
$ws.Range("H64").Value = 11112660
$ws.Range("I64").Value = 27778692
$ws.Range("J64").Value = 1971.5555
$ws.Range("K64").Value = 27778692
$ws.Range("L64").Value = 1971.5555
$ws.Range("M64").Value = -27778467
$ws.Range("N64").Value = -2421.5555
$ws.Range("H67").Value = 11112660
$ws.Range("I67").Value = 27778692
$ws.Range("J67").Value = 1971.5555
$ws.Range("K67").Value = 27778692
$ws.Range("L67").Value = 1971.5555
$ws.Range("M67").Value = -27777912
$ws.Range("N67").Value = -3531.5555
$ws.Range("H134").Value = 4723.815
$ws.Range("I134").Value = 1779.7428
$ws.Range("J134").Value = 10147.105
$ws.Range("K134").Value = 5339.2284
$ws.Range("L134").Value = 30441.315
$ws.Range("M134").Value = -2804.2284
$ws.Range("N134").Value = -35511.315

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5632.0625
$ws.Range("I16").Value = 3633.2856
$ws.Range("K16").Value = 3633.2856
$ws.Range("M16").Value = -3346.2856
$ws.Range("H58").Value = 11116380
$ws.Range("I58").Value = 21740410
$ws.Range("K58").Value = 21740410
$ws.Range("M58").Value = -21740207
$ws.Range("H94").Value = 626.625
$ws.Range("I94").Value = 505.63635
$ws.Range("J94").Value = 729
$ws.Range("K94").Value = 505.63635
$ws.Range("L94").Value = 729
$ws.Range("M94").Value = -54.63634999999999
$ws.Range("N94").Value = -1631
$ws.Range("H113").Value = 5632.0625
$ws.Range("I113").Value = 3633.2856
$ws.Range("K113").Value = 3633.2856
$ws.Range("M113").Value = -1463.2856
$ws.Range("H136").Value = 11116380
$ws.Range("I136").Value = 21740410
$ws.Range("K136").Value = 65221230
$ws.Range("M136").Value = -65218680
$ws.Range("H137").Value = 53615.91
$ws.Range("J137").Value = 51477.6
$ws.Range("L137").Value = 51477.6
$ws.Range("N137").Value = -61677.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 650
$ws.Range("I22").Value = 650
$ws.Range("K22").Value = 1950
$ws.Range("M22").Value = -1781
$ws.Range("H27").Value = 650
$ws.Range("I27").Value = 650
$ws.Range("K27").Value = 1950
$ws.Range("M27").Value = -1848
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2685
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -1908
$ws.Range("H75").Value = 66669172
$ws.Range("J75").Value = 37039820
$ws.Range("L75").Value = 111119460
$ws.Range("N75").Value = -111121456
$ws.Range("H78").Value = 66669172
$ws.Range("J78").Value = 37039820
$ws.Range("L78").Value = 333358380
$ws.Range("N78").Value = -333368364
$ws.Range("H128").Value = 219998.8
$ws.Range("I128").Value = 219998.8
$ws.Range("K128").Value = 659996.3999999999
$ws.Range("M128").Value = -655016.3999999999
$ws.Range("H134").Value = 5026.9
$ws.Range("I134").Value = 4724.3335
$ws.Range("K134").Value = 14173.0005
$ws.Range("M134").Value = -9103.000499999998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9322.579
$ws.Range("I70").Value = 8012.25
$ws.Range("J70").Value = 10275.546
$ws.Range("K70").Value = 8012.25
$ws.Range("L70").Value = 10275.546
$ws.Range("M70").Value = -7742.25
$ws.Range("N70").Value = -10815.546
$ws.Range("H73").Value = 9322.579
$ws.Range("I73").Value = 8012.25
$ws.Range("J73").Value = 10275.546
$ws.Range("K73").Value = 8012.25
$ws.Range("L73").Value = 10275.546
$ws.Range("M73").Value = -7076.25
$ws.Range("N73").Value = -12147.546

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 35000
$ws.Range("J33").Value = 35000
$ws.Range("L33").Value = 35000
$ws.Range("N33").Value = -35580
$ws.Range("H40").Value = 5258.9565
$ws.Range("I40").Value = 4006.6
$ws.Range("K40").Value = 4006.6
$ws.Range("M40").Value = -3870.6
$ws.Range("H93").Value = 7227.9287
$ws.Range("I93").Value = 6928.5713
$ws.Range("J93").Value = 7527.2856
$ws.Range("K93").Value = 6928.5713
$ws.Range("L93").Value = 7527.2856
$ws.Range("M93").Value = -5680.5713
$ws.Range("N93").Value = -10023.2856

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1448.7142
$ws.Range("I113").Value = 563.3333
$ws.Range("J113").Value = 3042.4
$ws.Range("K113").Value = 1689.9999
$ws.Range("L113").Value = 9127.200000000001
$ws.Range("M113").Value = 480.0001
$ws.Range("N113").Value = -13467.2
$ws.Range("H126").Value = 2665.2917
$ws.Range("I126").Value = 1610.5
$ws.Range("J126").Value = 4774.875
$ws.Range("K126").Value = 4831.5
$ws.Range("L126").Value = 14324.625
$ws.Range("M126").Value = -2361.5
$ws.Range("N126").Value = -19264.625
$ws.Range("H132").Value = 20850266
$ws.Range("I132").Value = 31257572
$ws.Range("J132").Value = 35654.25
$ws.Range("K132").Value = 93772716
$ws.Range("L132").Value = 106962.75
$ws.Range("M132").Value = -93770186
$ws.Range("N132").Value = -112022.75
$ws.Range("H136").Value = 25646062
$ws.Range("I136").Value = 55556468
$ws.Range("K136").Value = 166669404
$ws.Range("M136").Value = -166666854
